$wb = $excel.ActiveWorkbook

# --- TS_Defs sheet: the "p,t" tag in Q6 is narrowed down to just "t" ---
$wsDefs = $wb.Worksheets.Item("TS_Defs")
$wsDefs.Range("Q6").Value = "t"

# --- process map sheet: two new mapping rows appended (old_new -> */new and ep*/old) ---
$wsProc = $wb.Worksheets.Item("process map")
$wsProc.Range("A24").Value = "old_new"
$wsProc.Range("B24").Value = "*"
$wsProc.Range("C24").Value = "new"
$wsProc.Range("A25").Value = "old_new"
$wsProc.Range("B25").Value = "ep*"
$wsProc.Range("C25").Value = "old"

# --- Update saved selections to match where the author left the cursor ---
$wsDefs.Range("A6").Select()
$wsProc.Range("C26").Select()

# --- The "process map" tab was the active sheet when the workbook was last saved ---
$wsProc.Activate()
